$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" to "UK related" ---
$wsUkRelated = $wb.Worksheets.Item(1)
$wsUkRelated.Name = "UK related"

# --- Other countries sheet: update frozen-pane scroll position and selection ---
$wsOther = $wb.Worksheets.Item(2)
$wsOther.Activate()
[void]$wsOther.Range("M2").Select()

# Row heights changed (list restructured -> extra wrapped line for most rows)
$wsOther.Rows.Item(20).RowHeight = 85
$wsOther.Rows.Item(22).RowHeight = 85
$wsOther.Rows.Item(23).RowHeight = 102
$wsOther.Rows.Item(24).RowHeight = 85
$wsOther.Rows.Item(25).RowHeight = 85
$wsOther.Rows.Item(26).RowHeight = 85
$wsOther.Rows.Item(27).RowHeight = 85
$wsOther.Rows.Item(28).RowHeight = 85
$wsOther.Rows.Item(29).RowHeight = 102
$wsOther.Rows.Item(30).RowHeight = 85
$wsOther.Rows.Item(31).RowHeight = 68
$wsOther.Rows.Item(32).RowHeight = 102
$wsOther.Rows.Item(33).RowHeight = 102
$wsOther.Rows.Item(34).RowHeight = 102
$wsOther.Rows.Item(35).RowHeight = 102

# --- UK related sheet: update selection ---
$wsUkRelated.Activate()
[void]$wsUkRelated.Range("C4").Select()

# --- UK sheet: selection stays F41, just re-activate to normalize the view ---
$wsUk = $wb.Worksheets.Item(3)
$wsUk.Activate()
[void]$wsUk.Range("F41").Select()

# Leave "Other countries" as the active/visible tab, matching tabSelected in target.
$wsOther.Activate()
